$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 14 de Junio de 2020 a las 13:49"

# Row 4: 'Estados Unidos' -> 'Estados Unidos'
$ws.Range("B4").Value = 2142453
$ws.Range("C4").Value = 229
$ws.Range("E4").Value = 1170814
$ws.Range("G4").Value = 6
$ws.Range("H4").Value = 117533

# Row 7: 'India' -> 'India'
$ws.Range("B7").Value = 322647
$ws.Range("C7").Value = 1021
$ws.Range("D7").Value = 162709
$ws.Range("E7").Value = 150733
$ws.Range("G7").Value = 6
$ws.Range("H7").Value = 9205

# Row 26: 'Bielorrusia' -> 'Bielorrusia'
$ws.Range("B26").Value = 53973
$ws.Range("C26").Value = 732
$ws.Range("D26").Value = 30103
$ws.Range("E26").Value = 23562
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = 308

# Row 38: 'Suiza' -> 'Suiza'
$ws.Range("B38").Value = 31117
$ws.Range("C38").Value = 23
$ws.Range("E38").Value = 379

# Row 61: 'Ghana' -> 'Ghana'
$ws.Range("B61").Value = 11422
$ws.Range("C61").Value = 304
$ws.Range("D61").Value = 4156
$ws.Range("E61").Value = 7215
$ws.Range("G61").Value = 3
$ws.Range("H61").Value = 51

# Row 74: 'Nepal' -> 'Nepal'
$ws.Range("B74").Value = 5760
$ws.Range("C74").Value = 425
$ws.Range("D74").Value = 974
$ws.Range("E74").Value = 4767
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = 19

# Row 75: 'Senegal' -> 'Senegal'
$ws.Range("B75").Value = 5090
$ws.Range("C75").Value = 94
$ws.Range("D75").Value = 3344
$ws.Range("E75").Value = 1686

# Row 76: 'Uzbekistan' -> 'Uzbekistan'
$ws.Range("D76").Value = 3910
$ws.Range("E76").Value = 1065

# Row 79: 'Consejo Danes para los Refugiados' -> 'Consejo Danes para los Refugiados'
$ws.Range("B79").Value = 4778
$ws.Range("C79").Value = 54
$ws.Range("D79").Value = 600
$ws.Range("E79").Value = 4071
$ws.Range("G79").Value = 1
$ws.Range("H79").Value = 107

# Row 98: 'Croacia' -> 'Croacia'
$ws.Range("B98").Value = 2252
$ws.Range("C98").Value = 1
$ws.Range("E98").Value = 11

# Row 103: 'Sri Lanka' -> 'Sri Lanka'
$ws.Range("B103").Value = 1889
$ws.Range("C103").Value = 5
$ws.Range("E103").Value = 591

# Row 116: 'Libano' -> 'Libano'
$ws.Range("B116").Value = 1446
$ws.Range("C116").Value = 4
$ws.Range("E116").Value = 546

# Row 119: 'Paraguay' -> 'Madagascar'
$ws.Range("A119").Value = "Madagascar"
$ws.Range("B119").Value = 1272
$ws.Range("C119").Value = 20
$ws.Range("D119").Value = 367
$ws.Range("E119").Value = 895
$ws.Range("H119").Value = 10

# Row 120: 'Madagascar' -> 'Paraguay'
$ws.Range("A120").Value = "Paraguay"
$ws.Range("B120").Value = 1261
$ws.Range("D120").Value = 647
$ws.Range("E120").Value = 603
$ws.Range("H120").Value = 11

# Row 141: 'Malta' -> 'Malta'
$ws.Range("B141").Value = 649
$ws.Range("C141").Value = 3
$ws.Range("D141").Value = 603
$ws.Range("E141").Value = 37

# Row 210: 'Seychelles' -> 'Montserrat'
$ws.Range("A210").Value = "Montserrat"
$ws.Range("D210").Value = 10
$ws.Range("H210").Value = 1

# Row 211: 'Montserrat' -> 'Seychelles'
$ws.Range("A211").Value = "Seychelles"
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0

# Row 213: 'Islas Virgenes Britanicas' -> 'Papua Nueva Guinea'
$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("D213").Value = 8
$ws.Range("H213").Value = 0

# Row 214: 'Papua Nueva Guinea' -> 'Islas Virgenes Britanicas'
$ws.Range("A214").Value = "Islas Virgenes Britanicas"
$ws.Range("D214").Value = 7
$ws.Range("H214").Value = 1
